$d = $word.ActiveDocument

# Target the first paragraph (the **ID__...** marker paragraph)
$para = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt space, matching
# the target w:pBdr with w:space="5" on each edge.
$para.Range.Borders.DistanceFromTop = 5
$para.Range.Borders.DistanceFromLeft = 5
$para.Range.Borders.DistanceFromBottom = 5
$para.Range.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$para.Range.ParagraphFormat.LeftIndent = 11.25

# Remove the trailing space run that followed the marker text (before the
# paragraph mark).
$start = $para.Range.End - 2
$end = $para.Range.End - 1
$trailing = $d.Range($start, $end)
$trailing.Text = ""

# Update the ID marker text.
$d.Content.Find.Execute("**ID__AFFARS_5332_topic_3__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5332_104__ID**", 2)
